$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42
$ws.Cells.Item($row, 1).Value = "05/01/2026 05:37:39"
$ws.Cells.Item($row, 2).Value = "05/01 05:00"
$ws.Cells.Item($row, 3).Value = "Metrópoles"
$ws.Cells.Item($row, 4).Value = "Esquerda deve eleger um senador em SP, aposta cúpula do MDB"
$ws.Cells.Item($row, 5).Value = "https://www.metropoles.com/colunas/igor-gadelha/esquerda-deve-eleger-um-senador-em-sp-aposta-cupula-do-mdb"
$ws.Cells.Item($row, 6).Value = "senado"
$ws.Cells.Item($row, 7).Value = "Caciques do MDB apostam que esquerda deve eleger ao menos um &lt;b&gt;senado&lt;/b&gt;r em SP nas eleições de 2026, quando duas vagas estarão em disputa"
